$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values stay as text (avoid numeric auto-conversion),
# matching the original inlineStr cell type, then restore default style.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '62.418.09'
$ws.Range("E2").Value = '  -1.96%  '
$ws.Range("D3").Value = '3.161.87'
$ws.Range("E3").Value = '  -3.71%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '585.96'
$ws.Range("E5").Value = '  -3.14%  '
$ws.Range("D6").Value = '134.82'
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.157.32'
$ws.Range("E8").Value = '  -3.85%  '
$ws.Range("E9").Value = '  -2.26%  '
$ws.Range("E10").Value = '  -5.40%  '
$ws.Range("D11").Value = '5.25'
$ws.Range("E11").Value = '  -2.96%  '
$ws.Range("E12").Value = '  -3.21%  '
$ws.Range("D13").Value = '0.0000233'
$ws.Range("E13").Value = '  -4.87%  '
$ws.Range("E14").Value = '  -3.49%  '
$ws.Range("D15").Value = '3.682.78'
$ws.Range("E15").Value = '  -3.74%  '
$ws.Range("E16").Value = '  -2.03%  '
$ws.Range("D17").Value = '3.162.58'
$ws.Range("E17").Value = '  -3.69%  '
$ws.Range("D18").Value = '62.396.02'
$ws.Range("E18").Value = '  -2.10%  '
$ws.Range("E19").Value = '  -4.57%  '
$ws.Range("D20").Value = '454.34'
$ws.Range("E20").Value = '  -5.12%  '
$ws.Range("E21").Value = '  -0.98%  '
$ws.Range("E22").Value = '  -3.86%  '
$ws.Range("E23").Value = '  -5.04%  '
$ws.Range("D24").Value = '83.52'
$ws.Range("E24").Value = '  -0.75%  '
$ws.Range("D25").Value = '13.19'
$ws.Range("E25").Value = '  -2.45%  '
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("E28").Value = '  -3.27%  '
$ws.Range("E29").Value = '  -6.47%  '
$ws.Range("D30").Value = '7.70'
$ws.Range("E30").Value = '  -4.85%  '
$ws.Range("E31").Value = '  -7.07%  '
$ws.Range("D32").Value = '27.13'
$ws.Range("E32").Value = '  -5.38%  '
$ws.Range("E33").Value = '  -1.59%  '
$ws.Range("D34").Value = '2.38'
$ws.Range("E34").Value = '  -5.84%  '
$ws.Range("E35").Value = '  -6.46%  '
$ws.Range("D36").Value = '5.91'
$ws.Range("E36").Value = '  -0.63%  '
$ws.Range("D37").Value = '51.05'
$ws.Range("E37").Value = '  -4.45%  '
$ws.Range("D38").Value = '0.0₃0695'
$ws.Range("E38").Value = '  -5.70%  '
$ws.Range("E39").Value = '  -3.65%  '
$ws.Range("D40").Value = '2.73'
$ws.Range("E40").Value = '  +0.26%  '
$ws.Range("D41").Value = '396.51'
$ws.Range("E41").Value = '  -7.21%  '
$ws.Range("E42").Value = '  -4.07%  '
$ws.Range("E43").Value = '  -0.48%  '
$ws.Range("D44").Value = '2.791.99'
$ws.Range("E44").Value = '  -8.56%  '
$ws.Range("E45").Value = '  -5.58%  '
$ws.Range("D47").Value = '2.13'
$ws.Range("E47").Value = '  -2.35%  '
$ws.Range("E48").Value = '  +2.39%  '
$ws.Range("D49").Value = '125.19'
$ws.Range("E49").Value = '  +0.35%  '
$ws.Range("D50").Value = '25.22'
$ws.Range("E50").Value = '  -3.52%  '
$ws.Range("E51").Value = '  -3.88%  '

# Restore the default (unstyled) cell style for column D so no stray
# number-format styling is introduced that was not in the original file.
$ws.Range("D2:D51").Style = "Normal"
